$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the order codes in column C (Mã ĐH) from MU# to TU#
$ws.Range("C2").Value = "TU1"
$ws.Range("C3").Value = "TU2"
$ws.Range("C4").Value = "TU3"

# Update the selected cell on the sheet from C5 to F10
$ws.Range("F10").Select()
